$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 504, shifting existing rows 504-619 down to 505-620.
$ws.Rows.Item(504).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A504").Value = 11
$ws.Range("B504").Value = "Vega Monumental Concepción"
$ws.Range("C504").Value = "Bíobío"
$ws.Range("D504").Value = 44889
$ws.Range("E504").Value = 8
$ws.Range("F504").Value = 100112004
$ws.Range("G504").Value = "Cebolla"
$ws.Range("H504").Value = "Sin especificar"
$ws.Range("I504").Value = "1a nueva(o)"
$ws.Range("J504").Value = 200
$ws.Range("K504").Value = 12000
$ws.Range("L504").Value = 13000
$ws.Range("M504").Value = 12400
$ws.Range("N504").Value = "`$/malla 18 kilos"
$ws.Range("O504").Value = "Región Metropolitana"
$ws.Range("P504").Value = 689
$ws.Range("Q504").Value = 18
$ws.Range("R504").Value = "Hortaliza"

# Match the date number format already used in column D (style index 2 in styles.xml).
$ws.Range("D504").NumberFormat = $ws.Range("D505").NumberFormat
